# Scheduled-runner data refresh: updates Universalis market-price derived
# columns (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])
# for a handful of leve rows across the job sheets (H:N = columns 8-14).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 346.4762
$ws.Range("I33").Value = 332.75
$ws.Range("J33").Value = 390.4
$ws.Range("K33").Value = 332.75
$ws.Range("L33").Value = 390.4
$ws.Range("M33").Value = -103.75
$ws.Range("N33").Value = -848.4

$ws.Range("H58").Value = 1573.7142
$ws.Range("I58").Value = 1399.75
$ws.Range("J58").Value = 1805.6666
$ws.Range("K58").Value = 4199.25
$ws.Range("L58").Value = 5416.9998
$ws.Range("M58").Value = -4049.25
$ws.Range("N58").Value = -5716.9998

$ws.Range("H76").Value = 3400.1765
$ws.Range("I76").Value = 3376.8462
$ws.Range("J76").Value = 3476
$ws.Range("K76").Value = 3376.8462
$ws.Range("L76").Value = 3476
$ws.Range("M76").Value = -3061.8462
$ws.Range("N76").Value = -4106

$ws.Range("H79").Value = 3400.1765
$ws.Range("I79").Value = 3376.8462
$ws.Range("J79").Value = 3476
$ws.Range("K79").Value = 3376.8462
$ws.Range("L79").Value = 3476
$ws.Range("M79").Value = -2284.8462
$ws.Range("N79").Value = -5660

$ws.Range("H131").Value = 1942.2051
$ws.Range("I131").Value = 511.25
$ws.Range("J131").Value = 2105.743
$ws.Range("K131").Value = 1533.75
$ws.Range("L131").Value = 6317.228999999999
$ws.Range("M131").Value = 3506.25
$ws.Range("N131").Value = -16397.229

$ws.Range("H132").Value = 3128.1052
$ws.Range("I132").Value = 2452.6956
$ws.Range("J132").Value = 4163.7334
$ws.Range("K132").Value = 7358.0868
$ws.Range("L132").Value = 12491.2002
$ws.Range("M132").Value = -4828.0868
$ws.Range("N132").Value = -17551.2002

$ws.Range("H138").Value = 4446144.5
$ws.Range("I138").Value = 1300.875
$ws.Range("J138").Value = 9525966
$ws.Range("K138").Value = 3902.625
$ws.Range("L138").Value = 28577898
$ws.Range("M138").Value = 1237.375
$ws.Range("N138").Value = -28588178

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20366.861
$ws.Range("I32").Value = 19013.072
$ws.Range("J32").Value = 25556.389
$ws.Range("K32").Value = 19013.072
$ws.Range("L32").Value = 25556.389
$ws.Range("M32").Value = -18726.072
$ws.Range("N32").Value = -26130.389

$ws.Range("H61").Value = 27834984
$ws.Range("I61").Value = 32291578
$ws.Range("J61").Value = 204100
$ws.Range("K61").Value = 32291578
$ws.Range("L61").Value = 204100
$ws.Range("M61").Value = -32291366
$ws.Range("N61").Value = -204524

$ws.Range("H74").Value = 10481030
$ws.Range("I74").Value = 12860366
$ws.Range("K74").Value = 12860366
$ws.Range("M74").Value = -12859492

$ws.Range("H77").Value = 10481030
$ws.Range("I77").Value = 12860366
$ws.Range("K77").Value = 64301830
$ws.Range("M77").Value = -64297462

$ws.Range("H132").Value = 98099.17999999999
$ws.Range("I132").Value = 61952
$ws.Range("J132").Value = 220999.6
$ws.Range("K132").Value = 185856
$ws.Range("L132").Value = 662998.8
$ws.Range("M132").Value = -183326
$ws.Range("N132").Value = -668058.8

$ws.Range("H136").Value = 27834984
$ws.Range("I136").Value = 32291578
$ws.Range("J136").Value = 204100
$ws.Range("K136").Value = 96874734
$ws.Range("L136").Value = 612300
$ws.Range("M136").Value = -96872184
$ws.Range("N136").Value = -617400

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3720.32
$ws.Range("I31").Value = 4301.5
$ws.Range("K31").Value = 4301.5
$ws.Range("M31").Value = -4006.5

$ws.Range("H34").Value = 3720.32
$ws.Range("I34").Value = 4301.5
$ws.Range("K34").Value = 4301.5
$ws.Range("M34").Value = -4099.5

$ws.Range("H132").Value = 27179.426
$ws.Range("I132").Value = 1821.742
$ws.Range("J132").Value = 114522.555
$ws.Range("K132").Value = 5465.226
$ws.Range("L132").Value = 343567.665
$ws.Range("M132").Value = -2935.226
$ws.Range("N132").Value = -348627.665

$ws.Range("H138").Value = 39686.668
$ws.Range("J138").Value = 39686.668
$ws.Range("L138").Value = 39686.668
$ws.Range("N138").Value = -49966.668

$ws.Range("H139").Value = 51628.285
$ws.Range("J139").Value = 51628.285
$ws.Range("L139").Value = 51628.285
$ws.Range("N139").Value = -61908.285

$ws.Range("H141").Value = 61500
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 61500
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 61500
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -71860

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 1933.3334
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 1933.3334
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 5800.0002
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -6302.0002

$ws.Range("H75").Value = 2166
$ws.Range("I75").Value = 1632.25
$ws.Range("J75").Value = 2699.75
$ws.Range("K75").Value = 4896.75
$ws.Range("L75").Value = 8099.25
$ws.Range("M75").Value = -3898.75
$ws.Range("N75").Value = -10095.25

$ws.Range("H78").Value = 2166
$ws.Range("I78").Value = 1632.25
$ws.Range("J78").Value = 2699.75
$ws.Range("K78").Value = 14690.25
$ws.Range("L78").Value = 24297.75
$ws.Range("M78").Value = -9698.25
$ws.Range("N78").Value = -34281.75

$ws.Range("H103").Value = 2862.9656
$ws.Range("I103").Value = 439.25
$ws.Range("J103").Value = 3786.2856
$ws.Range("K103").Value = 1317.75
$ws.Range("L103").Value = 11358.8568
$ws.Range("M103").Value = -438.75
$ws.Range("N103").Value = -13116.8568

$ws.Range("H107").Value = 649.69446
$ws.Range("I107").Value = 535.2727
$ws.Range("J107").Value = 829.5
$ws.Range("K107").Value = 1605.8181
$ws.Range("L107").Value = 2488.5
$ws.Range("M107").Value = 314.1819
$ws.Range("N107").Value = -6328.5

$ws.Range("H121").Value = 65861816
$ws.Range("I121").Value = 1652
$ws.Range("J121").Value = 80179240
$ws.Range("K121").Value = 4956
$ws.Range("L121").Value = 240537720
$ws.Range("M121").Value = -3646
$ws.Range("N121").Value = -240540340

$ws.Range("H131").Value = 841.55554
$ws.Range("I131").Value = 410.81818
$ws.Range("J131").Value = 1031.08
$ws.Range("K131").Value = 1232.45454
$ws.Range("L131").Value = 3093.24
$ws.Range("M131").Value = 3807.54546
$ws.Range("N131").Value = -13173.24

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H80").Value = 3993.7273
$ws.Range("I80").Value = 2005
$ws.Range("K80").Value = 2005
$ws.Range("M80").Value = -1007

$ws.Range("H83").Value = 3993.7273
$ws.Range("I83").Value = 2005
$ws.Range("K83").Value = 10025
$ws.Range("M83").Value = -5033

$ws.Range("H113").Value = 2181.1177
$ws.Range("I113").Value = 1300
$ws.Range("J113").Value = 2964.3333
$ws.Range("K113").Value = 1300
$ws.Range("L113").Value = 2964.3333
$ws.Range("M113").Value = 870
$ws.Range("N113").Value = -7304.3333

$ws.Range("H132").Value = 68194.07000000001
$ws.Range("I132").Value = 41313.36
$ws.Range("J132").Value = 202597.6
$ws.Range("K132").Value = 123940.08
$ws.Range("L132").Value = 607792.8
$ws.Range("M132").Value = -121410.08
$ws.Range("N132").Value = -612852.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 461.09525
$ws.Range("I93").Value = 461.09525
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 461.09525
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 786.9047499999999
$ws.Range("N93").ClearContents()

$ws.Range("H132").Value = 86325.5
$ws.Range("I132").Value = 2238.5
$ws.Range("K132").Value = 6715.5
$ws.Range("M132").Value = -4185.5

$ws.Range("H136").Value = 46053.652
$ws.Range("I136").Value = 39585.297
$ws.Range("J136").Value = 55245.527
$ws.Range("K136").Value = 118755.891
$ws.Range("L136").Value = 165736.581
$ws.Range("M136").Value = -116205.891
$ws.Range("N136").Value = -170836.581

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1919.5
$ws.Range("I126").Value = 1862.875
$ws.Range("J126").Value = 1995
$ws.Range("K126").Value = 5588.625
$ws.Range("L126").Value = 5985
$ws.Range("M126").Value = -3118.625
$ws.Range("N126").Value = -10925

$ws.Range("H132").Value = 42343.55
$ws.Range("I132").Value = 26208.6
$ws.Range("J132").Value = 114054.445
$ws.Range("K132").Value = 78625.79999999999
$ws.Range("L132").Value = 342163.335
$ws.Range("M132").Value = -76095.79999999999
$ws.Range("N132").Value = -347223.335

$ws.Range("H136").Value = 50632.73
$ws.Range("I136").Value = 35498.69
$ws.Range("J136").Value = 87206.664
$ws.Range("K136").Value = 106496.07
$ws.Range("L136").Value = 261619.992
$ws.Range("M136").Value = -103946.07
$ws.Range("N136").Value = -266719.992
